# Demo for article progress bar
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: "Research" task now has a "Time spent" value of 0.5
$ws.Range("D8").Value = 0.5

# Row 9: Task renamed from "Codepen design" to "Design"; Time spent = 3
$ws.Range("B9").Value = "Design"
$ws.Range("D9").Value = 3

# Row 10: "Time spent" value of 0
$ws.Range("D10").Value = 0

# Update the active cell selection to D11 (cosmetic, matches author's cursor position)
$ws.Range("D11").Select()
